$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet and update the "through" label for the in-progress December 2021 month
$ws.Name = "Through 2021-12-26"
$ws.Range("B1").Value = "December 2021 (through December 26)"

# Updated counts (new data point added for 2022-01-03)
$ws.Range("N3").Value = 4     # Englewood, December 2020
$ws.Range("AX4").Value = 6    # North Lawndale, December 2017
$ws.Range("AX5").Value = 4    # West Pullman, December 2017
$ws.Range("N6").Value = 11    # Garfield Park, December 2020
$ws.Range("BJ6").Value = 7    # Garfield Park, December 2016
$ws.Range("B7").Value = 8     # Austin, December 2021 (through Dec 26)
$ws.Range("Z7").Value = 7     # Austin, December 2019
$ws.Range("AL7").Value = 13   # Austin, December 2018
$ws.Range("BV7").Value = 4    # Austin, December 2015
$ws.Range("AX9").Value = 2    # Grand Crossing, December 2017
$ws.Range("BJ9").Value = 1    # Grand Crossing, December 2016
$ws.Range("BV11").Value = 5   # Humboldt Park, December 2015
$ws.Range("N13").Value = 5    # Roseland, December 2020
$ws.Range("Z21").Value = 2    # Wicker Park, December 2019
$ws.Range("Z43").Value = 1    # Gage Park, December 2019
$ws.Range("AL43").Value = 1   # Gage Park, December 2018
$ws.Range("AX52").Value = 1   # Rush & Division, December 2017
$ws.Range("B53").Value = 2    # Hyde Park, December 2021 (through Dec 26)
$ws.Range("BJ82").Value = 3   # Logan Square, December 2016
$ws.Range("B93").Value = 4    # River North, December 2021 (through Dec 26)
